$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections on existing cells ---
$ws.Range("D8").Value = "I can view that my appointment is confirmed by me or not and when I should visit the clinic"

# --- Fill in missing "As a / an" cell on row 17 (Doctor) ---
$ws.Range("B17").Value = "Doctor"

# --- New rows for additional Doctor user stories ---
$ws.Range("B18").Value = "Doctor"
$ws.Range("C18").Value = "confirm if vaccination took place"
$ws.Range("D18").Value = "other process members can see if vaccination took place"

$ws.Range("B19").Value = "Doctor"
$ws.Range("C19").Value = "view appointment details"
$ws.Range("D19").Value = "I can clarify date and time of appointment"

$ws.Range("B20").Value = "Doctor"
$ws.Range("C20").Value = "view patient info"
$ws.Range("D20").Value = "I can get know which vaccination patient need"

# --- Widen the Clinic "Notes" entry to include the new type/vaccination fields ---
$ws.Range("E15").Value = "address, name, type, vacination types"

# --- Apply same wrap-text style used by columns C/D/E to the new rows ---
$ws.Range("C18:D18").WrapText = $true
$ws.Range("C19:D19").WrapText = $true
$ws.Range("C20:D20").WrapText = $true

# --- Row heights matching the surrounding wrapped rows (row 18 wraps to two
#     lines just like rows 16/17; 19/20 fit on one line at the default height) ---
$ws.Rows.Item(18).RowHeight = 30

# --- Minor column auto-fit widening that Excel performed after the new,
#     slightly-wider text was entered into columns C and E ---
$ws.Columns.Item(3).ColumnWidth = 40.33
$ws.Columns.Item(5).ColumnWidth = 34.33

# --- Update view / selection to match post-edit state ---
$ws.Range("B16:B20").Select()
$excel.ActiveWindow.ScrollRow = 13
